$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.942.33'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").Value = '1.633.12'
$ws.Range("E3").Value = '  -0.66%  '
$ws.Range("D4").Value = "'0.997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.42%  '
$ws.Range("D5").Value = "'211.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.84%  '
$ws.Range("E6").Value = '  -0.46%  '
$ws.Range("D7").Value = "'0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.45%  '
$ws.Range("D8").Value = "'23.41"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.49%  '
$ws.Range("E9").Value = '  -2.18%  '
$ws.Range("D10").Value = "'0.0614"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.42%  '
$ws.Range("E11").Value = '  +0.39%  '
$ws.Range("D12").Value = '1.863.62'
$ws.Range("E12").Value = '  -0.72%  '
$ws.Range("D13").Value = '1.622.25'
$ws.Range("E13").Value = '  -1.39%  '
$ws.Range("D14").Value = "'4.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.36%  '
$ws.Range("E15").Value = '  -2.04%  '
$ws.Range("D16").Value = "'65.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.56%  '
$ws.Range("D17").Value = '27.932.85'
$ws.Range("E17").Value = '  +0.13%  '
$ws.Range("D18").Value = "'232.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.74%  '
$ws.Range("E19").Value = '  +0.07%  '
$ws.Range("E20").Value = '  -0.99%  '
$ws.Range("D21").Value = "'0.995"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.47%  '
$ws.Range("D22").Value = "'10.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.06%  '
$ws.Range("D23").Value = "'4.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.96%  '
$ws.Range("E24").Value = '  -3.35%  '
$ws.Range("D25").Value = "'154.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.22%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("E27").Value = '  -0.84%  '
$ws.Range("D28").Value = "'15.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.43%  '
$ws.Range("D29").Value = "'0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.41%  '
$ws.Range("E30").Value = '  -1.10%  '
$ws.Range("D31").Value = "'0.0482"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.66%  '
$ws.Range("E32").Value = '  +1.98%  '
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("D34").Value = '1.406.04'
$ws.Range("E34").Value = '  -1.52%  '
$ws.Range("E35").Value = '  -0.19%  '
$ws.Range("D36").Value = "'1.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +8.42%  '
$ws.Range("E37").Value = '  +0.39%  '
$ws.Range("E38").Value = '  +1.63%  '
$ws.Range("D39").Value = "'0.558"
$ws.Range("D39").Style = "Normal"
$ws.Range("E40").Value = '  -1.58%  '
$ws.Range("E41").Value = '  -1.43%  '
$ws.Range("D42").Value = "'0.995"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.48%  '
$ws.Range("D43").Value = "'67.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.37%  '
$ws.Range("E44").Value = '  +2.16%  '
$ws.Range("E45").Value = '  +0.44%  '
$ws.Range("E46").Value = '  -0.56%  '
$ws.Range("D47").Value = '1.775.28'
$ws.Range("E47").Value = '  -0.57%  '
$ws.Range("D48").Value = "'88.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.04%  '
$ws.Range("E49").Value = '  -3.52%  '
$ws.Range("D50").Value = "'0.100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.63%  '
$ws.Range("E51").Value = '  -0.03%  '
